$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values that look numeric (e.g. "306.47", "1.001").
# Excel would otherwise auto-convert these text labels into real numbers,
# so force Text format on the column before writing, then restore the
# original (default) style afterwards so no stray number format sticks.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.197.59"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.904.16"
$ws.Range("E3").Value = "  +0.69%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "306.47"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "0.5256"
$ws.Range("E7").Value = "  +1.32%  "
$ws.Range("D8").Value = "0.3777"
$ws.Range("E8").Value = "  +1.47%  "
$ws.Range("D9").Value = "0.07257"
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("D10").Value = "21.12"
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("D11").Value = "0.8993"
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("D12").Value = "0.08417"
$ws.Range("E12").Value = "  +10.34%  "
$ws.Range("D13").Value = "1.890.73"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("D14").Value = "94.75"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "5.267"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "0.000008617"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").Value = "14.55"
$ws.Range("E18").Value = "  +1.26%  "
$ws.Range("D19").Value = "0.9998"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").Value = "27.226.84"
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D21").Value = "5.058"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "2.137.61"
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "6.433"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").Value = "146.85"
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("D26").Value = "2.276"
$ws.Range("E26").Value = "  +5.51%  "
$ws.Range("D27").Value = "1.751"
$ws.Range("E27").Value = "  -2.14%  "
$ws.Range("D28").Value = "18.16"
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("D29").Value = "114.74"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").Value = "4.927"
$ws.Range("E30").Value = "  -1.36%  "
$ws.Range("D31").Value = "4.802"
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("D33").Value = "0.8084"
$ws.Range("E33").Value = "  +6.40%  "
$ws.Range("D34").Value = "0.05062"
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("D35").Value = "1.235"
$ws.Range("E35").Value = "  +3.01%  "
$ws.Range("D36").Value = "2.949"
$ws.Range("E36").Value = "  -2.78%  "
$ws.Range("D37").Value = "3.354"
$ws.Range("E37").Value = "  +2.37%  "
$ws.Range("D38").Value = "2.606"
$ws.Range("E38").Value = "  +1.45%  "
$ws.Range("D39").Value = "0.5711"
$ws.Range("E39").Value = "  +1.20%  "
$ws.Range("D40").Value = "0.01988"
$ws.Range("E40").Value = "  -0.47%  "
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").Value = "6.649"
$ws.Range("E42").Value = "  +0.78%  "
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("D44").Value = "117.69"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("D46").Value = "0.4845"
$ws.Range("E46").Value = "  +0.46%  "
$ws.Range("D47").Value = "0.9995"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("D49").Value = "1.613"
$ws.Range("E49").Value = "  +2.16%  "
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("D51").Value = "63.60"
$ws.Range("E51").Value = "  +0.06%  "

$dRange.Style = "Normal"
